$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: _old -> _FV2310, _new -> _FV2404
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# Turn the used range into a native Excel Table ("Table1")
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U80"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1, freeze, select bottom-left pane)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

